$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the single-line mailing address into "street" / "city, state zip"
#    lines, and add a trailing blank paragraph after it.
$d.Content.Find.Execute("1028 Oak Tree Drive, San Jose CA 95129", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "1028 Oak Tree Drive^pSan Jose, CA 95129^p", 2)

# 3. Remove the two blank paragraphs that used to sit directly under
#    "Board of Directors" (one No Spacing, one Title-styled).
$anchor = $d.Content.Find
$anchor.Execute("Board of Directors", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = $anchor.Parent.Paragraphs(1)
$p.Next().Next().Range.Delete()
$p.Next().Range.Delete()
